$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-10 Tuesday", "2026-02-11 Wednesday"),
    @("702÷2=351, 0", "857÷6=142, 5"),
    @("701÷2=350, 1", "328÷3=109, 1"),
    @("645÷4=161, 1", "605÷9=67, 2"),
    @("513÷7=73, 2", "995÷9=110, 5"),
    @("924÷6=154, 0", "880÷2=440, 0"),
    @("501÷5=100, 1", "489÷5=97, 4"),
    @("822÷8=102, 6", "892÷4=223, 0"),
    @("402÷7=57, 3", "360÷2=180, 0"),
    @("501÷9=55, 6", "552÷9=61, 3"),
    @("299÷3=99, 2", "781÷9=86, 7"),
    @("222÷7=31, 5", "398÷6=66, 2"),
    @("196÷6=32, 4", "358÷2=179, 0"),
    @("797÷6=132, 5", "963÷7=137, 4"),
    @("490÷4=122, 2", "649÷7=92, 5"),
    @("636÷2=318, 0", "484÷2=242, 0"),
    @("445÷9=49, 4", "916÷7=130, 6"),
    @("637÷2=318, 1", "377÷8=47, 1"),
    @("381÷2=190, 1", "785÷7=112, 1"),
    @("249÷3=83, 0", "547÷6=91, 1"),
    @("975÷6=162, 3", "783÷2=391, 1"),
    @("115÷6=19, 1", "883÷9=98, 1"),
    @("233÷6=38, 5", "117÷3=39, 0"),
    @("353÷4=88, 1", "701÷3=233, 2"),
    @("833÷7=119, 0", "546÷8=68, 2"),
    @("516÷3=172, 0", "637÷5=127, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
